$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "investigaciones": update description for id 128 (row 34)
# ---------------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("investigaciones")
$wsInv.Range("C34").Value = "Caracterización de los hacedores y hacedoras de oficios artesanales que realizan aprovechamiento económico del espacio público en Bogotá. A partir de los aprendizajes del piloto realizado en 2024 y del marco normativo definido por el Decreto 315 de 2024 y la Resolución 500 de 2025, el estudio ajusta su objetivo, redefine la población de interés y actualiza el instrumento de recolección. Los resultados buscan aportar una comprensión más precisa de las condiciones de trabajo, los oficios, la cadena de valor y la relación de estas prácticas con el espacio público, como insumo para la gestión y formulación de políticas públicas."

# ---------------------------------------------------------------------------
# Sheet "productos": fill in previously-empty product rows, fix titles & add URLs
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("productos")

# investigacion_id 99
$wsProd.Range("B548").Value = "Informe final"
$wsProd.Range("C548").Value = "Informe final Caracterización de las ESALES culturales de Bogotá 2025`n"
$wsProd.Range("D548").Value = "Sí"
$wsProd.Range("E548").Value = "https://drive.google.com/file/d/1V2LcYxrrvtQrwyVv-PxQ9jHBnJp-Q5lX/view?usp=drive_link"

$wsProd.Range("B549").Value = "Instrumento recolección"
$wsProd.Range("C549").Value = "Formulario ESALES 2025"
$wsProd.Range("D549").Value = "Sí"
$wsProd.Range("E549").Value = "https://drive.google.com/file/d/1MzaFzKFri8-Y3gbOaECJEpN3dTcbHzKV/view?usp=drive_link"

# investigacion_id 101
$wsProd.Range("C561").Value = "Formulario EBC - PACCP Final Indicadores Formación"
$wsProd.Range("D561").Value = "Sí"
$wsProd.Range("E561").Value = "https://docs.google.com/spreadsheets/d/1FdjJSejVqDjcXO3k5xNzxpCU2Gv_V3er/edit?usp=sharing&ouid=105090632649587320414&rtpof=true&sd=true"

$wsProd.Range("B562").Value = "Base de datos"
$wsProd.Range("B563").Value = "Informe cuantitativo"

# investigacion_id 102
$wsProd.Range("C567").Value = "Base anonimizada Ecosistema Cultural 24/7 Asistentes"
$wsProd.Range("E567").Value = "https://docs.google.com/spreadsheets/d/1robWu-87K25n70k06owJuRggNwMVnXLF/edit?usp=drive_link&ouid=114639277514087565011&rtpof=true&sd=true"

$wsProd.Range("C568").Value = "Base anonimizada Ecosistema Cultural 24/7 Oferta "
$wsProd.Range("E568").Value = "https://docs.google.com/spreadsheets/d/1JTi5z-x1PgmRrGThTV8H0w86uOk3pDdd/edit?usp=drive_link&ouid=114639277514087565011&rtpof=true&sd=true"

$wsProd.Range("C569").Value = "Resultados Asistentes"
$wsProd.Range("E569").Value = "https://docs.google.com/spreadsheets/d/1Haze9wH0-YlSDY1bHrWDfXJjMwc1LD7J/edit?usp=drive_link&ouid=114639277514087565011&rtpof=true&sd=true"

$wsProd.Range("C570").Value = "Resultados Oferta"
$wsProd.Range("E570").Value = "https://docs.google.com/spreadsheets/d/1eL-kvZajGI8undcq5WGwmKI7cSji0iyJ/edit?usp=drive_link&ouid=114639277514087565011&rtpof=true&sd=true"

$wsProd.Range("C571").Value = "FichaTécnica Asistentes"
$wsProd.Range("E571").Value = "https://drive.google.com/file/d/1BPmvMkVgSMtofKkz65LdTUdby-e0nI66/view?usp=drive_link"

$wsProd.Range("C572").Value = "Ficha Técnica Oferta"
$wsProd.Range("E572").Value = "https://drive.google.com/file/d/121wFHZbwRLw_NGTkQ8UjMg0rsYtFQ3yl/view?usp=drive_link"

$wsProd.Range("E573").Value = "https://drive.google.com/file/d/1QIL_a1eOSFhjpKgln0XgVpFeVcFhCG1n/view?usp=drive_link"

# investigacion_id 109
$wsProd.Range("B610").Value = "Informe final"
$wsProd.Range("C610").Value = "Informe Festival Centro 2025"
$wsProd.Range("D610").Value = "Sí"
$wsProd.Range("E610").Value = "https://drive.google.com/file/d/1OFk4HXgTpXjngSXZwFWsEeZxJvkpUpdI/view?usp=drive_link"

$wsProd.Range("B611").Value = "Informe cuantitativo"
$wsProd.Range("C611").Value = "Tablas de salida Festival Centro"
$wsProd.Range("D611").Value = "Sí"
$wsProd.Range("E611").Value = "https://docs.google.com/spreadsheets/d/1e2ga-N1piD18x0W9kwcv0lMJDrROL_dW/edit?usp=drive_link&ouid=114639277514087565011&rtpof=true&sd=true"

# investigacion_id 118
$wsProd.Range("B664").Value = "Instrumento recolección"
$wsProd.Range("C664").Value = "Formulario EBC - PACCP Final"
$wsProd.Range("D664").Value = "Sí"
$wsProd.Range("E664").Value = "https://docs.google.com/spreadsheets/d/1qsclq1ogNVWTc6Sh7swGgV-k7DZr-Wnj/edit?usp=drive_link&ouid=105090632649587320414&rtpof=true&sd=true"

$wsProd.Range("B665").Value = "Base de datos"
$wsProd.Range("B666").Value = "Informe cuantitativo"

# ---------------------------------------------------------------------------
# Sheet "hallazgos": insert new findings for investigaciones 99, 102 and 109
# ---------------------------------------------------------------------------
$wsHal = $wb.Worksheets.Item("hallazgos")

# --- Insert 4 rows at 275 for investigacion_id 99 --------------------------
$wsHal.Range("A275:A278").EntireRow.Insert()

$wsHal.Range("A275").Value = 99
$wsHal.Range("B275").Value = 1
$wsHal.Range("C275").Value = "Áreas de actuación y alcance territorial: alta concentración urbana y fuerte orientación formativa–escénica`n"
$wsHal.Range("D275").Value = "El diagnóstico evidencia que las ESALES culturales en Bogotá desarrollan mayoritariamente su actividad en entornos urbanos y con una concentración territorial en localidades centrales como Teusaquillo (8,2%), Chapinero (7,7%), La Candelaria (6,4%) y Santa Fe (6,4%). Aunque una proporción relevante de entidades reporta alcance nacional y metropolitano (41,6% en cada uno de estos niveles), el anclaje operativo y organizativo sigue siendo predominantemente local y urbano, como lo confirma el hecho de que el 76,1% desarrolla sus actividades en ámbito urbano y el 97,9% tiene su ciudad de residencia en Bogotá. `n"
$wsHal.Range("E275").Value = "76,10%"
$wsHal.Range("F275").Value = "ESALES que desarrollan sus actividades en el ámbito urbano"

$wsHal.Range("A276").Value = 99
$wsHal.Range("B276").Value = 2
$wsHal.Range("C276").Value = "Fortalezas estructurales: organizaciones consolidadas pero con alta concentración de funciones"
$wsHal.Range("D276").Value = "El sector se caracteriza por una base de entidades con trayectorias medias y largas, dado que el 30,3% tiene entre 10 y 20 años de existencia y el 28,2% supera los 20 años de constitución jurídica. No obstante, esta fortaleza convive con una alta concentración de roles en una misma persona, particularmente en la figura del representante legal, quien en el 100% de los casos cumple simultáneamente funciones de dirección o gerencia y, también en el 100%, asume múltiples roles dentro de la organización. `n"
$wsHal.Range("E276").Value = "28,20%"
$wsHal.Range("F276").Value = "ESALESs con más de 20 años de constitución jurídica"

$wsHal.Range("A277").Value = 99
$wsHal.Range("B277").Value = 3
$wsHal.Range("C277").Value = "Financiamiento: dependencia del recurso público con avances relevantes en ingresos propios`n"
$wsHal.Range("D277").Value = "En materia de sostenibilidad económica, el diagnóstico muestra como principal fortaleza la capacidad de acceso a recursos públicos, utilizados por el 68,3% de las ESALES y que representan, en promedio, el 39,4% del total de su financiamiento. A ello se suma un avance en la venta de productos y/o servicios, reportada por el 57,0% de las entidades y que aporta en promedio el 38,7% de sus ingresos, lo que indica procesos incipientes pero relevantes de autogeneración de recursos. La empresa privada participa como fuente de financiamiento en el 40,9% de los casos, con un peso promedio del 15,0%. Sin embargo, persiste una baja diversificación financiera, evidenciada en la escasa participación de la cooperación internacional (9,2% de las entidades; 2,9% del financiamiento promedio) y de los créditos comerciales o internos (7,0%; 4,0% del promedio).`nEsta estructura financiera refuerza la estabilidad de corto plazo, pero mantiene al sector expuesto a riesgos asociados a la variabilidad de la inversión pública y a las limitaciones en capacidades comerciales y de negociación, en un contexto donde el 68,3% de las entidades reporta ingresos anuales inferiores a `$100 millones.`n"
$wsHal.Range("E277").Value = "68,30%"
$wsHal.Range("F277").Value = "ESALES que acceden a recursos públicos"

$wsHal.Range("A278").Value = 99
$wsHal.Range("B278").Value = 4
$wsHal.Range("C278").Value = "Perspectivas de mejora: consenso en la necesidad de fortalecer sostenibilidad y capacidades estratégicas`n"
$wsHal.Range("D278").Value = "Las ESALES culturales presentan una notable convergencia en sus perspectivas de futuro, al identificar como retos prioritarios la mejora de la sostenibilidad financiera (97,9% la considera importante o muy importante), el fortalecimiento organizacional (94,4%) y el aumento del número de beneficiarios o personas impactadas (94,4%). Este consenso se complementa con un interés creciente en la incorporación de herramientas digitales y tecnológicas, señaladas como innovación futura por el 42,3% de las entidades, así como en el fortalecimiento de capacidades relacionadas con la contratación pública (57,8%) y las alianzas estratégicas (57,0%) como temas prioritarios de capacitación. Al mismo tiempo, la baja definición de innovaciones futuras en el 44,4% de las organizaciones (Ns/Nr) sugiere la necesidad de procesos de acompañamiento estratégico que permitan traducir estas aspiraciones en planes concretos de desarrollo, diversificación territorial y ampliación del alcance cultural.`n"
$wsHal.Range("E278").Value = ""
$wsHal.Range("F278").Value = ""

# --- Insert 5 rows at 284 for investigacion_id 102 (after id 100 / id 103 block shift) ---
$wsHal.Range("A284:A288").EntireRow.Insert()

$wsHal.Range("A284").Value = 102
$wsHal.Range("B284").Value = 1
$wsHal.Range("C284").Value = "El ecosistema cultural nocturno en Bogotá presenta condiciones para una ampliación horaria gradual"
$wsHal.Range("D284").Value = "Si bien una proporción de establecimientos opera en la franja nocturna temprana y la ciudadanía manifiesta alta disposición a asistir a actividades en la noche, la operación sostenida en horarios extendidos y de madrugada es limitada, tanto en la práctica como en la disposición declarada."
$wsHal.Range("E284").Value = ""
$wsHal.Range("F284").Value = ""

$wsHal.Range("A285").Value = 102
$wsHal.Range("B285").Value = 2
$wsHal.Range("C285").Value = "Existe una desalineación entre la disposición de la demanda y las capacidades operativas de la oferta"
$wsHal.Range("D285").Value = "Mientras la ciudadanía expresa altos niveles de interés y disposición a asistir con frecuencia regular a eventos nocturnos, la oferta presenta limitaciones asociadas a protocolos, articulación entre actores y gestión del riesgo, lo que condiciona su capacidad de responder a dicha demanda."
$wsHal.Range("E285").Value = ""
$wsHal.Range("F285").Value = ""

$wsHal.Range("A286").Value = 102
$wsHal.Range("B286").Value = 3
$wsHal.Range("C286").Value = "La estrategia Ecosistema Cultural 24/7 presenta una baja apropiación desde la ciudadanía, pese a una valoración mayoritariamente positiva del enfoque"
$wsHal.Range("D286").Value = "El bajo nivel de conocimiento de la estrategia contrasta con la aceptación general de la ampliación de la oferta nocturna, lo que sugiere una brecha entre el diseño institucional y su posicionamiento público."
$wsHal.Range("E286").Value = ""
$wsHal.Range("F286").Value = ""

$wsHal.Range("A287").Value = 102
$wsHal.Range("B287").Value = 4
$wsHal.Range("C287").Value = "Las barreras para la operación y el consumo nocturno son principalmente de carácter urbano y estructural"
$wsHal.Range("D287").Value = "La seguridad y la movilidad emergen como los principales factores que condicionan tanto la decisión de los establecimientos para ampliar horarios como la disposición de la ciudadanía para asistir, por encima de restricciones normativas o culturales."
$wsHal.Range("E287").Value = ""
$wsHal.Range("F287").Value = ""

$wsHal.Range("A288").Value = 102
$wsHal.Range("B288").Value = 5
$wsHal.Range("C288").Value = "El ecosistema se compone mayoritariamente de establecimientos de pequeña y mediana escala, con formas de operación que tienden a desarrollarse de manera individual y localizada"
$wsHal.Range("D288").Value = "Este perfil favorece dinámicas de consumo de cercanía y experiencias de menor escala, especialmente en franjas horarias extendidas, aunque también plantea desafíos para la articulación entre actores y la construcción de ofertas nocturnas integrales, particularmente en ausencia de mecanismos de coordinación o acompañamiento institucional."
$wsHal.Range("E288").Value = ""
$wsHal.Range("F288").Value = ""

# --- Insert 4 rows at 298 for investigacion_id 109 (after id 105 block shift) ---
$wsHal.Range("A298:A301").EntireRow.Insert()

$wsHal.Range("A298").Value = 109
$wsHal.Range("B298").Value = 1
$wsHal.Range("C298").Value = "El Festival activa dinámicas económicas y de sosenibilidad cultural"
$wsHal.Range("D298").Value = "Aunque el apoyo a emprendimientos culturales, la compra de productos o el intercambio con melómanos no constituyen el principal motivo de asistencia en ninguno de los escenarios, sí aparecen de manera consistente en todos ellos como motivaciones complementarias. Este patrón sugiere que el Festival Centro no solo funciona como un espacio de circulación artística, sino también como un dispositivo que contribuye a la sostenibilidad económica del ecosistema cultural, al activar prácticas de consumo cultural, visibilización de agentes y circulación de bienes simbólicos asociados a la música y las artes.`nAsí, se podría seguir fortaleciendo de manera estratégica los componentes de circulación económica y visibilización de emprendimientos culturales, especialmente en aquellos escenarios con mayor afluencia y diversidad de públicos, sin desdibujar el eje artístico del Festival.`n"

$wsHal.Range("A299").Value = 109
$wsHal.Range("B299").Value = 2
$wsHal.Range("C299").Value = "El Festival no genera mayores impactos negativos en el espacio público"
$wsHal.Range("D299").Value = "En los cuatro escenarios analizados, la percepción mayoritaria de los asistentes indica que la realización del Festival no modifica sustancialmente problemáticas asociadas al espacio público, como el arrojo de basuras, el parqueo en zonas prohibidas o la contaminación auditiva y visual. `nPor otro lado, en algunos casos, como el Muelle de la FUGA y La Media Torta, se registra una mayor percepción de incremento en la presencia de vendedores informales o en el turismo; sin embargo, estos fenómenos coexisten con una valoración positiva del evento y no se asocian a un deterioro de la convivencia en la zona.`n"

$wsHal.Range("A300").Value = 109
$wsHal.Range("B300").Value = 3
$wsHal.Range("C300").Value = "El Festival articula los desplazamientos, recorridos y consumos culturales en el centro"
$wsHal.Range("D300").Value = "En todos los escenarios, una proporción mayoritaria de asistentes declaró haber visitado o tener previsto visitar otros espacios del centro antes o después del evento. Los recorridos se concentraron principalmente en equipamientos culturales, cafés, restaurantes, bares, teatros, salas de arte y museos, lo que evidencia una alta capacidad del Festival para articular la oferta cultural, gastronómica y comercial del centro de Bogotá. Este comportamiento refuerza el papel del Festival Centro como dinamizador territorial y como nodo de conexión entre la programación cultural pública y otras actividades económicas y simbólicas del área."

$wsHal.Range("A301").Value = 109
$wsHal.Range("B301").Value = 4
$wsHal.Range("C301").Value = "La oferta musical y artística es el principal motivo de asistencia`n"
$wsHal.Range("D301").Value = "En todos los escenarios, los motivos de asistencia se concentran de manera consistente en la presencia de los grupos y artistas, el interés por conocer nuevas propuestas musicales y el reconocimiento previo de algunos de los artistas participantes. El Festival Centro opera simultáneamente como un espacio de encuentro entre artistas y sus audiencias y como una plataforma para el descubrimiento de nuevas propuestas, lo que explica la coexistencia de públicos con trayectoria en el Festival y de personas que asisten por primera vez."
